$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update meandam (E), lowerdam (F), upperdam (G) values for the power_tower rows.
$ws.Range("E3").Value = 115532.8092832113
$ws.Range("F3").Value = 86649.60696240845
$ws.Range("G3").Value = 144416.0116040141

$ws.Range("E5").Value = 118784.4471689444
$ws.Range("F5").Value = 89088.3353767083
$ws.Range("G5").Value = 148480.5589611805

$ws.Range("E7").Value = 156793.3506980647
$ws.Range("F7").Value = 117595.0130235485
$ws.Range("G7").Value = 195991.6883725809

$ws.Range("E9").Value = 168686.7014348303
$ws.Range("F9").Value = 126515.0260761227
$ws.Range("G9").Value = 210858.3767935378

$ws.Range("E11").Value = 180933.3065996079
$ws.Range("F11").Value = 135699.9799497059
$ws.Range("G11").Value = 226166.6332495098

$ws.Range("E13").Value = 191392.426329325
$ws.Range("F13").Value = 143544.3197469937
$ws.Range("G13").Value = 239240.5329116562

$ws.Range("E15").Value = 204245.0672493127
$ws.Range("F15").Value = 153183.8004369845
$ws.Range("G15").Value = 255306.3340616409

$ws.Range("E17").Value = 219330.4923569412
$ws.Range("F17").Value = 164497.8692677059
$ws.Range("G17").Value = 274163.1154461765

$ws.Range("E19").Value = 230803.7639313215
$ws.Range("F19").Value = 173102.8229484911
$ws.Range("G19").Value = 288504.7049141519

$ws.Range("E21").Value = 243434.4708011228
$ws.Range("F21").Value = 182575.8531008421
$ws.Range("G21").Value = 304293.0885014035

$wb.Save()
